$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update selection to AB35 (was T23)
$ws.Range("AB35").Select()

$ws.Range("AE5").Value = 5
$ws.Range("AE6").Value = 5
$ws.Range("AE10").Value = 5
$ws.Range("AE11").Value = 5
# Row 12: full data row
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 15
$ws.Range("G12").Value = 5
$ws.Range("H12").Value = 5
$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 5
$ws.Range("K12").Value = 20
$ws.Range("L12").Value = 5
$ws.Range("M12").Value = 5
$ws.Range("N12").Value = 5
$ws.Range("O12").Value = 15
$ws.Range("P12").Value = 5
$ws.Range("Q12").Value = 5
$ws.Range("R12").Value = 5
$ws.Range("S12").Value = 5
$ws.Range("T12").Value = 20
$ws.Range("U12").Value = 5
$ws.Range("V12").Value = 5
$ws.Range("W12").Value = 5
$ws.Range("X12").Value = 15
$ws.Range("Y12").Value = 5
$ws.Range("Z12").Value = 5
$ws.Range("AA12").Value = 5
$ws.Range("AB12").Value = 5
$ws.Range("AC12").Value = 20
$ws.Range("AD12").Value = 5
$ws.Range("AE12").Value = 5

# Row 13: full data row
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = 13
$ws.Range("G13").Value = 5
$ws.Range("H13").Value = 4
$ws.Range("I13").Value = 5
$ws.Range("J13").Value = 5
$ws.Range("K13").Value = 19
$ws.Range("L13").Value = 5
$ws.Range("M13").Value = 5
$ws.Range("N13").Value = 5
$ws.Range("O13").Value = 15
$ws.Range("P13").Value = 5
$ws.Range("Q13").Value = 5
$ws.Range("R13").Value = 5
$ws.Range("S13").Value = 5
$ws.Range("T13").Value = 20
$ws.Range("U13").Value = 5
$ws.Range("V13").Value = 5
$ws.Range("W13").Value = 5
$ws.Range("X13").Value = 15
$ws.Range("Y13").Value = 4
$ws.Range("Z13").Value = 4
$ws.Range("AA13").Value = 4
$ws.Range("AB13").Value = 4
$ws.Range("AC13").Value = 16
$ws.Range("AD13").Value = 5
$ws.Range("AE13").Value = 5

$ws.Range("AE15").Value = 5
$ws.Range("AE16").Value = 4
$ws.Range("AE17").Value = 5
$ws.Range("AE18").Value = 5
$ws.Range("AE19").Value = 5
$ws.Range("AE20").Value = 5
$ws.Range("AE21").Value = 5
$ws.Range("AE22").Value = 5
$ws.Range("AE23").Value = 5
$ws.Range("AE24").Value = 5
# Row 25: full data row
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 5
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 5
$ws.Range("I25").Value = 5
$ws.Range("J25").Value = 5
$ws.Range("K25").Value = 20
$ws.Range("L25").Value = 5
$ws.Range("M25").Value = 5
$ws.Range("N25").Value = 5
$ws.Range("O25").Value = 15
$ws.Range("P25").Value = 5
$ws.Range("Q25").Value = 5
$ws.Range("R25").Value = 5
$ws.Range("S25").Value = 5
$ws.Range("T25").Value = 20
$ws.Range("U25").Value = 5
$ws.Range("V25").Value = 5
$ws.Range("W25").Value = 5
$ws.Range("X25").Value = 15
$ws.Range("Y25").Value = 5
$ws.Range("Z25").Value = 5
$ws.Range("AA25").Value = 5
$ws.Range("AB25").Value = 5
$ws.Range("AC25").Value = 20
$ws.Range("AD25").Value = 5
$ws.Range("AE25").Value = 5

$ws.Range("AE27").Value = 5
$ws.Range("AE28").Value = 5
$ws.Range("AE29").Value = 4
$ws.Range("AE30").Value = 4
$ws.Range("AE31").Value = 5
$ws.Range("AE32").Value = 5
$ws.Range("AE33").Value = 5
$ws.Range("AE34").Value = 5
$ws.Range("AE35").Value = 5
$ws.Range("AE36").Value = 5
$ws.Range("AE37").Value = 4
$ws.Range("AE38").Value = 4
$ws.Range("AE39").Value = 5
$ws.Range("AE40").Value = 5
$ws.Range("AE41").Value = 5
$ws.Range("AE42").Value = 5
$ws.Range("AE43").Value = 5
$ws.Range("AE44").Value = 5
$ws.Range("AE45").Value = 5
$ws.Range("AE46").Value = 5
$ws.Range("AE47").Value = 5
$ws.Range("AE48").Value = 4
$ws.Range("AE49").Value = 4
